# Update the "Förändrad" (changed) date column (C) for rows 2-44
# from 2023-09-13 (serial 45182) to 2023-09-15 (serial 45184).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 44; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
